$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Diagonal of single-letter labels "a".."g" starting at D2, stepping one
# row/column at a time down to J8 (mirrors the existing A/B "Name"/letter
# table already on the sheet).
$ws.Range("D2").Value = "a"
$ws.Range("E3").Value = "b"
$ws.Range("F4").Value = "c"
$ws.Range("G5").Value = "d"
$ws.Range("H6").Value = "e"
$ws.Range("I7").Value = "f"
$ws.Range("J8").Value = "g"

# A second diagonal of numeric values two rows below the table, D10..I15.
$ws.Range("D10").Value = 11.1
$ws.Range("E11").Value = 22.2
$ws.Range("F12").Value = 33.3
$ws.Range("G13").Value = 44.4
$ws.Range("H14").Value = 55.5
$ws.Range("I15").Value = 66.6

# Leave the selection on the last cell touched, same as the source edit.
$ws.Range("I15").Select()
